$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Drop the trailing "size_min" parameter block (columns AJ:AS). The new
#    layout only has 35 columns (A:AI) instead of 45 (A:AS).
# ---------------------------------------------------------------------------
$ws.Range("AJ1:AS11").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2. Re-label the header row (row 2) with the v2 "coef" parameter names and
#    update the associated sample row (row 4) with the matching values.
# ---------------------------------------------------------------------------
$headers = @(
    "coef15.0","coef17.0","coef19.0","coef111.0","coef113.0",
    "coef25.0","coef27.0","coef29.0","coef211.0","coef213.0",
    "coef35.0","coef37.0","coef39.0","coef311.0","coef313.0",
    "coef45.0","coef47.0","coef49.0","coef411.0","coef413.0",
    "coef55.0","coef57.0","coef59.0","coef511.0","coef513.0",
    "coef80.0","coef81.0","coef82.0","coef83.0","coef84.0",
    "coef125.0","coef127.0","coef129.0","coef1211.0","coef1213.0"
)

$row4values = @(
    5,7,9,11,13,
    5,7,9,11,13,
    5,7,9,11,13,
    5,7,9,11,13,
    5,7,9,11,13,
    0,1,2,3,4,
    5,7,9,11,13
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(2, $col).Value = $headers[$i]
    $ws.Cells.Item(4, $col).Value = $row4values[$i]
}
